$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14: swap order of elements in the set literal shown in E14
$ws.Range("E14").Value = "{'Tuple[NoneType]', 'Tuple[None]'}"

# Row 15: change Tuple[None] -> Tuple[NoneType]
$ws.Range("E15").Value = "Tuple[NoneType]"

# Row 45: move "Scalpel Accuracy:" label from C45 to E45, clear old D45 value,
# and put the corrected accuracy value in F45 (matching the layout used in row 46)
$ws.Range("C45").Value = ""
$ws.Range("D45").Value = ""
$ws.Range("E45").Value = "Scalpel Accuracy:"
$ws.Range("F45").Value = 95.24000000000001

# Row 46: fix label wording
$ws.Range("E46").Value = "Accuracy vs PyType"
